# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column after "Week"
#  - shorten the Week labels (W01 -> W1, etc.)
#  - fill in the new Week_Start_Date values
#  - correct the MyForecast values
#  - store is_holiday_week as a boolean
# Also refresh the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old "ASIN" column (B) for Week_Start_Date.
$ws.Range("B1").EntireColumn.Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Format as text first so the ISO-like dates aren't auto-converted to date serials.
$ws.Range("B2:B17").NumberFormat = "@"

$weeks = @("W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9", "W10", "W11", "W12", "W13", "W14", "W15", "W16")

$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

$myForecast = @(22, 25, 24, 27, 25, 26, 24, 27, 25, 30, 22, 31, 22, 28, 27, 23)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]
    $ws.Cells.Item($row, 10).Value = $false
}

# Refresh the forecast totals on the Summary sheet (kept as text, matching
# the existing cells in that column).
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9:B12").NumberFormat = "@"
$summary.Cells.Item(9, 2).Value = "410"
$summary.Cells.Item(10, 2).Value = "201"
$summary.Cells.Item(11, 2).Value = "98"
$summary.Cells.Item(12, 2).Value = "31"
